$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 340.5
$ws.Range("J9").Value = 414.33334
$ws.Range("L9").Value = 414.33334
$ws.Range("N9").Value = -752.33334

$ws.Range("H38").Value = 1048.625
$ws.Range("I38").Value = 1048.625
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 3145.875
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -2773.875
$ws.Range("N38").ClearContents()

$ws.Range("H87").Value = 60000
$ws.Range("J87").Value = 60000
$ws.Range("L87").Value = 60000
$ws.Range("N87").Value = -62496

$ws.Range("H90").Value = 60000
$ws.Range("J90").Value = 60000
$ws.Range("L90").Value = 180000
$ws.Range("N90").Value = -192480

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1989530.1
$ws.Range("I32").Value = 2277578.5
$ws.Range("K32").Value = 2277578.5
$ws.Range("M32").Value = -2277291.5

$ws.Range("H45").Value = 4419.5
$ws.Range("I45").Value = 1995.6666
$ws.Range("J45").Value = 5873.8
$ws.Range("K45").Value = 1995.6666
$ws.Range("L45").Value = 5873.8
$ws.Range("M45").Value = -1618.6666
$ws.Range("N45").Value = -6627.8

$ws.Range("H61").Value = 5136.7817
$ws.Range("I61").Value = 2856.4358
$ws.Range("J61").Value = 10695.125
$ws.Range("K61").Value = 2856.4358
$ws.Range("L61").Value = 10695.125
$ws.Range("M61").Value = -2644.4358
$ws.Range("N61").Value = -11119.125

$ws.Range("H74").Value = 2764.6858
$ws.Range("I74").Value = 1642.0625
$ws.Range("J74").Value = 3710.0527
$ws.Range("K74").Value = 1642.0625
$ws.Range("L74").Value = 3710.0527
$ws.Range("M74").Value = -768.0625
$ws.Range("N74").Value = -5458.0527

$ws.Range("H77").Value = 2764.6858
$ws.Range("I77").Value = 1642.0625
$ws.Range("J77").Value = 3710.0527
$ws.Range("K77").Value = 8210.3125
$ws.Range("L77").Value = 18550.2635
$ws.Range("M77").Value = -3842.3125
$ws.Range("N77").Value = -27286.2635

$ws.Range("H110").Value = 13890093
$ws.Range("I110").Value = 1097.1111
$ws.Range("J110").Value = 55557080
$ws.Range("K110").Value = 1097.1111
$ws.Range("L110").Value = 55557080
$ws.Range("M110").Value = 947.8888999999999
$ws.Range("N110").Value = -55561170

$ws.Range("H132").Value = 3020.9736
$ws.Range("I132").Value = 1603.2623
$ws.Range("K132").Value = 4809.7869
$ws.Range("M132").Value = -2279.7869

$ws.Range("H136").Value = 5136.7817
$ws.Range("I136").Value = 2856.4358
$ws.Range("J136").Value = 10695.125
$ws.Range("K136").Value = 8569.307400000002
$ws.Range("L136").Value = 32085.375
$ws.Range("M136").Value = -6019.307400000002
$ws.Range("N136").Value = -37185.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 11113520
$ws.Range("I20").Value = 16668661
$ws.Range("J20").Value = 3239
$ws.Range("K20").Value = 16668661
$ws.Range("L20").Value = 3239
$ws.Range("M20").Value = -16668414
$ws.Range("N20").Value = -3733

$ws.Range("H86").Value = 55560212
$ws.Range("J86").Value = 90914510
$ws.Range("L86").Value = 90914510
$ws.Range("N86").Value = -90916756

$ws.Range("H89").Value = 55560212
$ws.Range("J89").Value = 90914510
$ws.Range("L89").Value = 454572550
$ws.Range("N89").Value = -454583782

$ws.Range("H105").Value = 42739.92
$ws.Range("I105").Value = 51502.418
$ws.Range("K105").Value = 51502.418
$ws.Range("M105").Value = -49755.418

$ws.Range("H134").Value = 4474.707
$ws.Range("I134").Value = 1690.5
$ws.Range("K134").Value = 5071.5
$ws.Range("M134").Value = -2536.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4752.619
$ws.Range("I16").Value = 3065.9092
$ws.Range("J16").Value = 6608
$ws.Range("K16").Value = 3065.9092
$ws.Range("L16").Value = 6608
$ws.Range("M16").Value = -2778.9092
$ws.Range("N16").Value = -7182

$ws.Range("H99").Value = 3826.0417
$ws.Range("J99").Value = 7138.4287
$ws.Range("L99").Value = 7138.4287
$ws.Range("N99").Value = -10134.4287

$ws.Range("H113").Value = 4752.619
$ws.Range("I113").Value = 3065.9092
$ws.Range("J113").Value = 6608
$ws.Range("K113").Value = 3065.9092
$ws.Range("L113").Value = 6608
$ws.Range("M113").Value = -895.9092000000001
$ws.Range("N113").Value = -10948

$ws.Range("H126").Value = 3826.0417
$ws.Range("J126").Value = 7138.4287
$ws.Range("L126").Value = 21415.2861
$ws.Range("N126").Value = -26355.2861

$ws.Range("H134").Value = 4142.9814
$ws.Range("I134").Value = 1174.5278
$ws.Range("J134").Value = 10079.889
$ws.Range("K134").Value = 3523.5834
$ws.Range("L134").Value = 30239.667
$ws.Range("M134").Value = -988.5834000000004
$ws.Range("N134").Value = -35309.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3443.6667
$ws.Range("I68").Value = 1500.2858
$ws.Range("J68").Value = 6164.4
$ws.Range("K68").Value = 4500.857400000001
$ws.Range("L68").Value = 18493.2
$ws.Range("M68").Value = -3689.857400000001
$ws.Range("N68").Value = -20115.2

$ws.Range("H71").Value = 3443.6667
$ws.Range("I71").Value = 1500.2858
$ws.Range("J71").Value = 6164.4
$ws.Range("K71").Value = 13502.5722
$ws.Range("L71").Value = 55479.6
$ws.Range("M71").Value = -9446.572200000001
$ws.Range("N71").Value = -63591.6

$ws.Range("H109").Value = 55557140
$ws.Range("I109").Value = 62501180
$ws.Range("J109").Value = 37039704
$ws.Range("K109").Value = 187503540
$ws.Range("L109").Value = 111119112
$ws.Range("M109").Value = -187502500
$ws.Range("N109").Value = -111121192

$ws.Range("H113").Value = 5786.2666
$ws.Range("I113").Value = 2193.7144
$ws.Range("K113").Value = 6581.1432
$ws.Range("M113").Value = -4411.1432

$ws.Range("H114").Value = 380.33334
$ws.Range("J114").Value = 1000
$ws.Range("L114").Value = 3000
$ws.Range("N114").Value = -9508

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 8417.966
$ws.Range("I132").Value = 4382.3076
$ws.Range("K132").Value = 39440.7684
$ws.Range("M132").Value = -36910.7684

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2306322.2
$ws.Range("I122").Value = 3107361.5
$ws.Range("J122").Value = 3334.75
$ws.Range("K122").Value = 9322084.5
$ws.Range("L122").Value = 10004.25
$ws.Range("M122").Value = -9319634.5
$ws.Range("N122").Value = -14904.25

$ws.Range("H126").Value = 3939.2632
$ws.Range("I126").Value = 4199.273
$ws.Range("K126").Value = 12597.819
$ws.Range("M126").Value = -10127.819

$ws.Range("H132").Value = 3594.8857
$ws.Range("I132").Value = 1522.2858
$ws.Range("K132").Value = 4566.857400000001
$ws.Range("M132").Value = -2036.857400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2174.2778
$ws.Range("I46").Value = 484
$ws.Range("J46").Value = 3019.4167
$ws.Range("K46").Value = 484
$ws.Range("L46").Value = 3019.4167
$ws.Range("M46").Value = -296
$ws.Range("N46").Value = -3395.4167

$ws.Range("H132").Value = 8777214
$ws.Range("I132").Value = 17243260
$ws.Range("K132").Value = 51729780
$ws.Range("M132").Value = -51727250

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 875.5
$ws.Range("J100").Value = 1017.125
$ws.Range("L100").Value = 2034.25
$ws.Range("N100").Value = -3116.25

$ws.Range("H122").Value = 4678.8096
$ws.Range("I122").Value = 4286.4136
$ws.Range("K122").Value = 12859.2408
$ws.Range("M122").Value = -10409.2408

$ws.Range("H126").Value = 2514
$ws.Range("I126").Value = 1509.1666
$ws.Range("K126").Value = 4527.4998
$ws.Range("M126").Value = -2057.4998

$ws.Range("H132").Value = 16680370
$ws.Range("I132").Value = 26322278
$ws.Range("K132").Value = 78966834
$ws.Range("M132").Value = -78964304
